# Fruta / hortaliza, semanal
# Re-shuffle the weekly rows: the data in columns D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# N (Unidad de comercializacion), O (Origen), P (Precio $/Kg) and
# Q (Kg o Unidades) move between rows as new weekly records replace the
# old ones. Row 7 is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the "before" values for the rows that get shuffled, keyed by
# their original row number, for columns D, J, K, L, M, N, O, P, Q.
$cols = @(4, 10, 11, 12, 13, 14, 15, 16, 17)   # D, J, K, L, M, N, O, P, Q

$snapshot = @{}
foreach ($r in @(2, 3, 4, 5, 6, 8, 9, 10)) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Mapping of destination row -> source row (data that now lands there)
$rowMap = @{
    2  = 5
    3  = 4
    4  = 8
    5  = 9
    6  = 3
    8  = 2
    9  = 10
    10 = 6
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
